# Scene.xlsx edit: finish quest swamp and flowsand
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Quest column (E) values for four scenes to add new quest entries.
# Order matters for shared-string table append order, so set cells in the
# same order the new strings appear in the target workbook.

# Row 18 - riverold (遗迹河滩): add swamp;2
$ws.Range("E18").Value = "river;2|fishpool;1|swamp;2"

# Row 21 - riverside (落潮小径): add swamp;2
$ws.Range("E21").Value = "mushroom;1|hiddeway;1|swamp;2"

# Row 19 - fogvalley (雾谷): add sandflow;2
$ws.Range("E19").Value = "portal;3|sandflow;2"

# Row 6  - orevalley (矿脉山脚): add sandflow;2
$ws.Range("E6").Value = "sandpile;1|stone;2|sandflow;2"

# Update the active cell selection to E7
$ws.Range("E7").Select()
